# Updates cryptos list values (Price and Volume(1h) columns) per upstream refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.298.68"
$ws.Range("E2").Value = "  -0.55%  "
$ws.Range("D3").Value = "1.883.06"
$ws.Range("E3").Value = "  -1.42%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "238.17"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.28%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.000"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.01%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4699"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.55%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2846"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.61%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06607"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.35%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "20.86"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +11.87%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07790"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.51%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "98.25"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.41%  "
$ws.Range("D13").Value = "1.882.06"
$ws.Range("E13").Value = "  -1.49%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.103"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.77%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6785"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.77%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "284.73"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +11.34%  "
$ws.Range("D17").Value = "30.320.04"
$ws.Range("E17").Value = "  -0.57%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.001"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.03%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.66"
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "5.416"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.76%  "
$ws.Range("D21").Value = "2.127.24"
$ws.Range("E21").Value = "  -1.39%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.000007294"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.08%  "
$ws.Range("E23").Value = "  -0.10%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.193"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.38%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.420"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.00%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "168.58"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.55%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "19.30"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.22%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.999"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.60%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.373"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.89%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.09731"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.82%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.411"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -7.45%  "
$ws.Range("E32").Value = "  -1.28%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.143"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.45%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.04683"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.36%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.7091"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.91%  "
$ws.Range("E36").Value = "  -0.25%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.718"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.62%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01881"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.48%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "6.676"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +6.62%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.526"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -3.24%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "72.31"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -3.54%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.982"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.21%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.8655"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.63%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.0000"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.02%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "103.35"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.90%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4206"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.57%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "992.38"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +8.02%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.296"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.88%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.216"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +5.28%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "34.04"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.83%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.1153"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.62%  "

Write-Output "Updated cryptos list"
